# Rename a few header cells across the workbook's sheets:
#   member_profiles!A1: "Member ID"          -> "MemberID"
#   member_profiles!D1: "SEBI Registration No" -> "SEBIRegistrationNo"
#   kmp!A1:                "Company"         -> "CompanyName"
#   authorized_personnel!A1: "Company"       -> "CompanyName"

$wb = $excel.ActiveWorkbook

$wsMembers = $wb.Worksheets.Item("member_profiles")
$wsKmp = $wb.Worksheets.Item("kmp")
$wsAuthorized = $wb.Worksheets.Item("authorized_personnel")

$wsKmp.Range("A1").Value = "CompanyName"
$wsAuthorized.Range("A1").Value = "CompanyName"
$wsMembers.Range("D1").Value = "SEBIRegistrationNo"
$wsMembers.Range("A1").Value = "MemberID"
